$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (Y1:AC1) need the same bold/border/centered style as the
# existing header row. Copy format from X1 (an existing styled header cell)
# before filling in the new labels, so no new cell styles are introduced.
$ws.Range("X1").Copy()
$ws.Range("Y1:AC1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row (row 1): set text labels for all columns B..AC
$ws.Range("B1").Value = 'Total Cost'
$ws.Range("C1").Value = 'CO2'
$ws.Range("D1").Value = 'crudeoil'
$ws.Range("E1").Value = 'crudeoilisFixed'
$ws.Range("F1").Value = 'hydrogen'
$ws.Range("G1").Value = 'hydrogenisFixed'
$ws.Range("H1").Value = 'biomass'
$ws.Range("I1").Value = 'biomassisFixed'
$ws.Range("J1").Value = 'CrOilopex'
$ws.Range("K1").Value = 'H2opex'
$ws.Range("L1").Value = 'BMopex'
$ws.Range("M1").Value = 'Kilometersopex'
$ws.Range("N1").Value = 'KilometersDemand'
$ws.Range("O1").Value = 'Refineryopex'
$ws.Range("P1").Value = 'RefineryTotalEff'
$ws.Range("Q1").Value = 'Refinery-gasoline'
$ws.Range("R1").Value = 'Refinery-km'
$ws.Range("S1").Value = 'MtGopex'
$ws.Range("T1").Value = 'MtGTotalEff'
$ws.Range("U1").Value = 'MtG-gasoline'
$ws.Range("V1").Value = 'B2gasopex'
$ws.Range("W1").Value = 'B2gasTotalEff'
$ws.Range("X1").Value = 'B2gas-gasoline'
$ws.Range("Y1").Value = 'Gtkmopex'
$ws.Range("Z1").Value = 'GtkmTotalEff'
$ws.Range("AA1").Value = 'Gtkm-km'
$ws.Range("AB1").Value = 'GasHubopex'
$ws.Range("AC1").Value = 'KmHubopex'

# Data rows 2..11: column A is the index (unchanged), B..AC are numeric values
# Row 2
$ws.Range("B2").Value = 145.8620424750275
$ws.Range("C2").Value = 24.12795879239154
$ws.Range("D2").Value = 329.616923393327
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 2816.520842045401
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0.014028679211575
$ws.Range("K2").Value = 0.05014624791182692
$ws.Range("L2").Value = 0.04877195807954137
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 979.0187249532061
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0.8870626299992788
$ws.Range("Q2").Value = 292.0984641025979
$ws.Range("R2").Value = 0.1461954274787777
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0.8329786587329684
$ws.Range("U2").Value = 2346.101753300429
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0.5283213260742617
$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 0.3710379989617706
$ws.Range("AA2").Value = 978.8725295257273
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = 0

# Row 3
$ws.Range("B3").Value = 75.18892136327152
$ws.Range("C3").Value = 40
$ws.Range("D3").Value = 390.2731115987876
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 2179.412189043648
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 285.8002057742187
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0.001467405111490405
$ws.Range("K3").Value = 0.03206914779836689
$ws.Range("L3").Value = 0.01653022252952235
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 989.9426343696258
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0.9280575578643128
$ws.Range("Q3").Value = 361.8337149396268
$ws.Range("R3").Value = 0.1810979554252386
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0.8212188984479827
$ws.Range("U3").Value = 1789.774477150531
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0.4624227707518719
$ws.Range("X3").Value = 132.1605230355692
$ws.Range("Y3").Value = 0
$ws.Range("Z3").Value = 0.4333895678046853
$ws.Range("AA3").Value = 989.7615364142006
$ws.Range("AB3").Value = 0
$ws.Range("AC3").Value = 0

# Row 4
$ws.Range("B4").Value = 166.0449509395155
$ws.Range("C4").Value = 40
$ws.Range("D4").Value = 333.5350612908001
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 2568.94433721117
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 389.6308378378359
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0.01396003344996548
$ws.Range("K4").Value = 0.05918397578040363
$ws.Range("L4").Value = 0.02399309804204856
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 989.1392756832992
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0.9421562939632433
$ws.Range("Q4").Value = 313.9279150952909
$ws.Range("R4").Value = 0.1571210786262717
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0.822789303155954
$ws.Range("U4").Value = 2113.699921060413
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0.5628996776315032
$ws.Range("X4").Value = 219.3230730142104
$ws.Range("Y4").Value = 0
$ws.Range("Z4").Value = 0.3736307126734053
$ws.Range("AA4").Value = 988.9821546046729
$ws.Range("AB4").Value = 0
$ws.Range("AC4").Value = 0

# Row 5
$ws.Range("B5").Value = 164.658234198897
$ws.Range("C5").Value = 40
$ws.Range("D5").Value = 291.0951240841666
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 2495.235677258024
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 467.2959229259749
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0.01248242949934429
$ws.Range("K5").Value = 0.06430596882595495
$ws.Range("L5").Value = 0.001211463939603356
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 1015.362602977342
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0.9277863761547354
$ws.Range("Q5").Value = 269.8040162000716
$ws.Range("R5").Value = 0.135037045145181
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0.8415663659009266
$ws.Range("U5").Value = 2099.906420976373
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 0.5122912175367222
$ws.Range("X5").Value = 239.3915973056942
$ws.Range("Y5").Value = 0
$ws.Range("Z5").Value = 0.3891099514372583
$ws.Range("AA5").Value = 1015.227565932197
$ws.Range("AB5").Value = 0
$ws.Range("AC5").Value = 0

# Row 6
$ws.Range("B6").Value = 171.9664156053566
$ws.Range("C6").Value = 40
$ws.Range("D6").Value = 291.909370923222
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 2363.717672703871
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 465.8058512105036
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0.02091012626972292
$ws.Range("K6").Value = 0.06589207348089593
$ws.Range("L6").Value = 0.02170924902695353
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 998.6705953756284
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0.9446726007928928
$ws.Range("Q6").Value = 275.4830258412316
$ws.Range("R6").Value = 0.1378793923129287
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0.8227869082623929
$ws.Range("U6").Value = 1944.835955929197
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 0.5370490679487279
$ws.Range("X6").Value = 250.1605982376648
$ws.Range("Y6").Value = 0
$ws.Range("Z6").Value = 0.4041857799852951
$ws.Range("AA6").Value = 998.5327159833155
$ws.Range("AB6").Value = 0
$ws.Range("AC6").Value = 0

# Row 7
$ws.Range("B7").Value = 159.3668576390237
$ws.Range("C7").Value = 40
$ws.Range("D7").Value = 279.3284676361202
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 2470.492016324709
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 488.8289042259
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0.00951712593234633
$ws.Range("K7").Value = 0.05658343406475675
$ws.Range("L7").Value = 0.03461237904828417
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 1000.218219448921
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0.9166462039983477
$ws.Range("Q7").Value = 255.7893341477976
$ws.Range("R7").Value = 0.1280226897636625
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 0.8355639882676629
$ws.Range("U7").Value = 2064.254162143694
$ws.Range("V7").Value = 0
$ws.Range("W7").Value = 0.5538568575680402
$ws.Range("X7").Value = 270.7412407829852
$ws.Range("Y7").Value = 0
$ws.Range("Z7").Value = 0.3860182524806991
$ws.Range("AA7").Value = 1000.090196759157
$ws.Range("AB7").Value = 0
$ws.Range("AC7").Value = 0

# Row 8
$ws.Range("B8").Value = 165.7288151529837
$ws.Range("C8").Value = 21.82561331024009
$ws.Range("D8").Value = 298.1641162601105
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 2578.506393756364
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0.046872349121371
$ws.Range("K8").Value = 0.05885312635548715
$ws.Range("L8").Value = 0.0408011388190223
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 996.0863748800409
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0.8968413436485037
$ws.Range("Q8").Value = 267.1385007478316
$ws.Range("R8").Value = 0.1337029533272431
$ws.Range("S8").Value = 0
$ws.Range("T8").Value = 0.8478755814941091
$ws.Range("U8").Value = 2186.252607992455
$ws.Range("V8").Value = 0
$ws.Range("W8").Value = 0.5356691682654375
$ws.Range("X8").Value = 0
$ws.Range("Y8").Value = 0
$ws.Range("Z8").Value = 0.4059494095248815
$ws.Range("AA8").Value = 995.9526719267137
$ws.Range("AB8").Value = 0
$ws.Range("AC8").Value = 0

# Row 9
$ws.Range("B9").Value = 95.16940342584414
$ws.Range("C9").Value = 40
$ws.Range("D9").Value = 264.8536679659116
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 2182.83436390863
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 515.3177876223818
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0.02606388889867657
$ws.Range("K9").Value = 0.03963769808703445
$ws.Range("L9").Value = 0.003383848570422905
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 1002.758793064754
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0.9334423896161448
$ws.Range("Q9").Value = 246.9784150839768
$ws.Range("R9").Value = 0.1236128203623508
$ws.Range("S9").Value = 0
$ws.Range("T9").Value = 0.8652878331417447
$ws.Range("U9").Value = 1888.780016853837
$ws.Range("V9").Value = 0
$ws.Range("W9").Value = 0.5779152870105407
$ws.Range("X9").Value = 297.8100271354256
$ws.Range("Y9").Value = 0
$ws.Range("Z9").Value = 0.4120020443666578
$ws.Range("AA9").Value = 1002.635180244392
$ws.Range("AB9").Value = 0
$ws.Range("AC9").Value = 0

# Row 10
$ws.Range("B10").Value = 131.5415287974047
$ws.Range("C10").Value = 40
$ws.Range("D10").Value = 296.4910604663992
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 2187.268089317909
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 457.4213593464895
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0.009351994581335716
$ws.Range("K10").Value = 0.05419213510761991
$ws.Range("L10").Value = 0.02237765680253914
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 993.2067954237351
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0.9119126474184127
$ws.Range("Q10").Value = 270.1035739379209
$ws.Range("R10").Value = 0.1351869739429034
$ws.Range("S10").Value = 0
$ws.Range("T10").Value = 0.8698232942414507
$ws.Range("U10").Value = 1902.536734839708
$ws.Range("V10").Value = 0
$ws.Range("W10").Value = 0.5450470295004306
$ws.Range("X10").Value = 249.3161531418532
$ws.Range("Y10").Value = 0
$ws.Range("Z10").Value = 0.410028678906453
$ws.Range("AA10").Value = 993.0716084497922
$ws.Range("AB10").Value = 0
$ws.Range("AC10").Value = 0

# Row 11
$ws.Range("B11").Value = 102.4001388278409
$ws.Range("C11").Value = 40
$ws.Range("D11").Value = 311.0205671850397
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 2300.901625518436
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 430.8323620513773
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0.02209198340787805
$ws.Range("K11").Value = 0.03756355071854477
$ws.Range("L11").Value = 0.02111968252851739
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 1004.433346476786
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0.9107731287250183
$ws.Range("Q11").Value = 282.9859058978754
$ws.Range("R11").Value = 0.1416345875364742
$ws.Range("S11").Value = 0
$ws.Range("T11").Value = 0.850738753548603
$ws.Range("U11").Value = 1957.466180931509
$ws.Range("V11").Value = 0
$ws.Range("W11").Value = 0.544095820514717
$ws.Range("X11").Value = 234.4140875346382
$ws.Range("Y11").Value = 0
$ws.Range("Z11").Value = 0.4057963708471334
$ws.Range("AA11").Value = 1004.291711889249
$ws.Range("AB11").Value = 0
$ws.Range("AC11").Value = 0

